$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Insert a new column B (GME part names); shifts old B,C,D -> C,D,E ----
$ws.Columns("B:B").Insert()

# ---- 2. Drop the stray placeholder cell left behind in F1 by the shift ----
$ws.Range("F1").Clear()

# ---- 3. Fill the new GME column B with part names (rows 1-14) ----
$ws.Range("B1").Value = "GME"
$ws.Range("B2").Value = "KLS MLW10G konektor"
$ws.Range("B3").Value = "pruzinovy pogo kontakt samice"
$ws.Range("B4").Value = "CONNFLY S2G06C, kolikova lista 2x3 pin"
$ws.Range("B5").Value = "dutinkova lista 8 pin"
$ws.Range("B6").Value = "CONNFLY S1G02C kolíková lišta 1x2 pin"
$ws.Range("B7").Value = "BL104G-V5,7 dutinková lišta 1x4pin"
$ws.Range("B8").Value = "dutinková lišta 1x6pin 2.54mm roztec"
$ws.Range("B9").Value = "KLS BL105G-V8.5 dutinková lišta 1x5 pin"
$ws.Range("B10").Value = "CONNFLY S1G04C kolíková lišta 1x4 pin"
$ws.Range("B11").Value = "CONNFLY S1G06C kolíková lišta 1x6"
$ws.Range("B12").Value = " KLS S1G08C kolíková lišta 1x8"
$ws.Range("B13").Value = " KLS S1G10C kolíková lišta 1x10"
$ws.Range("B14").Value = "distancni slopek plastovy 10mm"

# ---- 4. New row 14 quantity ----
$ws.Range("A14").Value = 6

# ---- 5. New G3 cell (Digikey URL for the pogo-pin part) ----
$ws.Range("G3").Value = "https://www.digikey.com/en/products/detail/mill-max-manufacturing-corp/4141-0-00-15-00-00-03-0/16341759"

# ---- 6. Rebuild hyperlinks so they track the shifted Mouser-number/URL columns ----
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E3"), "https://cz.mouser.com/ProductDetail/Mill-Max/4141-0-00-15-00-00-03-0?qs=Znm5pLBrcAKrk2sGDtRJwQ%3D%3D ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://cz.mouser.com/ProductDetail/Molex/70246-1004?qs=R7%2FZKp6KZ2by8%252BhtUWbo%252BQ%3D%3D ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://cz.mouser.com/ProductDetail/Amphenol-FCI/10129381-906002BLF?qs=DXv0QSHKF4zSKXKyBbBVMw%3D%3D ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E6"), "https://cz.mouser.com/ProductDetail/Samtec/TSW-102-07-T-S?qs=hzq9O0YpnZuE%252BS9CaLPVlQ%3D%3D ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E9"), "https://cz.mouser.com/ProductDetail/Wurth-Elektronik/61300511821?qs=iLbezkQI%252BsjtLywwvTh3KA%3D%3D ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E5"), "https://cz.mouser.com/ProductDetail/Gravitech/8Fx1L-254mm?qs=fkzBJ5HM%252BdAyuablm941Ag%3D%3D") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E10"), "https://cz.mouser.com/ProductDetail/Amphenol-Commercial-Products/G800W304018EU?qs=f9yNj16SXrKi9QS16pMdvA%3D%3D ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E8"), "https://cz.mouser.com/ProductDetail/Wurth-Elektronik/61300611821?qs=iLbezkQI%252BsgS21jF2eozhA%3D%3D") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E11"), "https://cz.mouser.com/ProductDetail/Amphenol-FCI/10129378-906001BLF?qs=0lQeLiL1qybdgsDy2F9d0g%3D%3D ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E12"), "https://cz.mouser.com/ProductDetail/Amphenol-FCI/10129378-908001BLF?qs=0lQeLiL1qybDtVkZc9V5DA%3D%3D ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E13"), "https://cz.mouser.com/ProductDetail/Amphenol-FCI/10129378-910001BLF?qs=0lQeLiL1qyaKTim8PTf2WA%3D%3D ") | Out-Null

# ---- 7. Re-fit the data columns to their (new) content, like Excel does on edit ----
$ws.Columns("A:E").AutoFit()
$ws.Columns("F:F").ColumnWidth = 28.75

# ---- 8. Restore the saved selection/active cell ----
$ws.Range("B20").Select()
